$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the separator row's formatting (black fill) before row 4 gets new data,
#     by copying it down to the new blank row 7 ---
$ws.Range("A4:E4").Copy()
$ws.Range("A7:G7").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Row 4 no longer holds the separator formatting - it now holds real data below
$ws.Range("A4:E4").Style = "Normal"

# --- Populate the new / changed cells ---
$ws.Range("C4").Value = "Tim"
$ws.Range("D4").Value = "Bob"
$ws.Range("F1").Value = "nullValue"
$ws.Range("B2").Value = "No"
$ws.Range("A4").Value = "AddCustomerError1"
$ws.Range("G1").Value = "Desc"
$ws.Range("G2").Value = "Adds Customer No Errors"
$ws.Range("G4").Value = "Mandatory Field Error PostCode pops up"
$ws.Range("G5").Value = "Mandatory Field Error LastName pops up"
$ws.Range("A5").Value = "AddCustomerError2"
$ws.Range("A6").Value = "AddCustomerError3"
$ws.Range("D6").Value = "Gary"
$ws.Range("G6").Value = "Mandatory Field Error FirstName pops up"

$ws.Range("B3").Value = "No"
$ws.Range("F2").Value = "No"
$ws.Range("F3").Value = "No"
$ws.Range("G3").Value = "Adds Customer No Errors"

$ws.Range("B4").Value = "Yes"
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = "PostCode"

$ws.Range("B5").Value = "Yes"
$ws.Range("C5").Value = "John"
$ws.Range("E5").Value = 1011
$ws.Range("F5").Value = "LastName"

$ws.Range("B6").Value = "Yes"
$ws.Range("E6").Value = 1011

# --- Column widths (best-fit-like widths for the newly widened columns) ---
$ws.Columns.Item(1).ColumnWidth = 16.8
$ws.Columns.Item(7).ColumnWidth = 36.5

# --- Selection / view state ---
$ws.Range("I17").Select()
